$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing B/C values with refreshed ECB long-term forecast figures ---
$ws.Range("B5").Value = 3.499141862354735
$ws.Range("B6").Value = 3.0204455992343306
$ws.Range("B7").Value = 3.5125897781563267
$ws.Range("B8").Value = 3.1574271404874361
$ws.Range("B9").Value = 2.9732649834018154
$ws.Range("B10").Value = 2.7950670292830289
$ws.Range("B11").Value = 3.2527451461202257
$ws.Range("B12").Value = 3.136648886594092
$ws.Range("B14").Value = 3.4457490558095092
$ws.Range("B19").Value = 4.1964485503483546
$ws.Range("B22").Value = 3.4874202741550695
$ws.Range("B23").Value = 2.7036394912638082
$ws.Range("B24").Value = 4.0865494512064702
$ws.Range("B26").Value = 3.8672670582059214
$ws.Range("B28").Value = 3.8507329289080485
$ws.Range("B29").Value = 4.264923010122315
$ws.Range("B30").Value = 3.804812490003763
$ws.Range("C30").Value = 3.1500000953674316
$ws.Range("C32").Value = 3.4000000953674316
$ws.Range("B42").Value = 3.1614793661370681
$ws.Range("B43").Value = 3.6511622567549966
$ws.Range("B45").Value = 3.9115942924127554
$ws.Range("B47").Value = 4.3402053283443838
$ws.Range("B48").Value = 3.8853768983537429
$ws.Range("B49").Value = 3.9346189130210516
$ws.Range("B51").Value = 4.3105899595757666
$ws.Range("B52").Value = 4.3959122317039006

# --- Fill in B56:C57 (previously blank) and append new monthly rows 58-64 ---
# Copy the date/number formatting down from row 55 for the newly appended rows
$ws.Range("A55:C55").Copy()
$ws.Range("A58:C64").PasteSpecial(-4122)

$ws.Range("B56").Value = 4.1494526746385363
$ws.Range("C56").Value = 3.5

$ws.Range("B57").Value = 4.0994317907359088
$ws.Range("C57").Value = 3

$ws.Range("A58").Value = 45536
$ws.Range("B58").Value = 3.9088409532940434
$ws.Range("C58").Value = 3

$ws.Range("A59").Value = 45566
$ws.Range("B59").Value = 4.5244216438247786
$ws.Range("C59").Value = 3.4000000953674316

$ws.Range("A60").Value = 45597
$ws.Range("B60").Value = 4.534444939525204
$ws.Range("C60").Value = 3.2000000476837158

$ws.Range("A61").Value = 45627
$ws.Range("B61").Value = 4.6021189435766479
$ws.Range("C61").Value = 3.2000000476837158

$ws.Range("A62").Value = 45658
$ws.Range("B62").Value = 4.3238947451761458
$ws.Range("C62").Value = 3

$ws.Range("A63").Value = 45689
$ws.Range("B63").Value = 4.7441090386901772
$ws.Range("C63").Value = 4

$ws.Range("A64").Value = 45717
$ws.Range("B64").Value = 4.5395244355094615
$ws.Range("C64").Value = 3.5

